$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Test_Cases")
$ws2 = $wb.Worksheets.Item("Test_Data")

# --- Update Run_Mode (column C) values on Test_Cases ---
$ws1.Range("C4").Value = "N"
$ws1.Range("C5").Value = "Y"
$ws1.Range("C6").Value = "N"
$ws1.Range("C8").Value = "Y"
$ws1.Range("C11").Value = "Y"

# --- Apply AutoFilter on Testing_Type (column D) = "Regression", which hides rows 3,4,6,7,9,10 ---
if ($ws1.AutoFilterMode) {
    $ws1.AutoFilterMode = $false
}
$ws1.Range("A1:E11").AutoFilter(4, @("Regression"), 7)

# --- Keep the _FilterDatabase defined name in sync with the new filter range ---
$wb.Names.Item("Test_Cases!_FilterDatabase").RefersTo = "=Test_Cases!`$A`$1:`$E`$11"

# --- View/selection state ---
# Set Test_Data's own selection/scroll first, then switch back to Test_Cases so it stays the active tab.
$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 10
$ws2.Range("A26").Select()

$ws1.Activate()
$ws1.Range("B5").Select()
